$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-10 from 2023-10-25 (45224)
# to 2023-11-03 (45233), keeping existing number formatting/style intact.
$newDate = Get-Date -Year 2023 -Month 11 -Day 3 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
